# Insert a new weekly price record for "Vega Monumental Concepción - Coliflor"
# at row 438, pushing the existing rows 438:487 down to 439:488.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 438 (shifts rows 438-487 down to 439-488).
$ws.Rows.Item(438).Insert()

# Populate the newly inserted row 438 with the new record's data.
$ws.Cells.Item(438, 1).Value = 11
$ws.Cells.Item(438, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(438, 3).Value = "Bíobío"
$ws.Cells.Item(438, 4).Value = 45212
$ws.Cells.Item(438, 5).Value = 8
$ws.Cells.Item(438, 6).Value = 100112008
$ws.Cells.Item(438, 7).Value = "Coliflor"
$ws.Cells.Item(438, 8).Value = "Sin especificar"
$ws.Cells.Item(438, 9).Value = "Primera"
$ws.Cells.Item(438, 10).Value = 2000
$ws.Cells.Item(438, 11).Value = 900
$ws.Cells.Item(438, 12).Value = 1000
$ws.Cells.Item(438, 13).Value = 950
$ws.Cells.Item(438, 14).Value = "`$/unidad"
$ws.Cells.Item(438, 15).Value = "Región Metropolitana"
$ws.Cells.Item(438, 16).Value = 950
$ws.Cells.Item(438, 17).Value = 1
$ws.Cells.Item(438, 18).Value = "Hortaliza"
